$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sensor data (columns C-H) for rows 2-31; timestamps (A) and label (B) follow a simple pattern.
$data = @(
    @(2, 2.669549942016602, -27.01569938659668, 1.220865488052368, 0.8970553824457185, -0.2614919603881221, -0.8384511420282582),
    @(3, -9.358198165893556, -20.71603393554688, -7.516577243804932, -0.4747457169627652, 0.1235830856455553, 2.941044032415826),
    @(4, -13.51923274993896, -8.136168479919434, -7.273739814758301, 0.6183352183740025, -0.2533354612730362, 0.0438180289910041),
    @(5, -7.558236122131348, -11.54725742340088, -2.366233348846436, 0.2323529118577299, 1.99577405316345, -1.971809049310575),
    @(6, -7.292168617248535, -5.43182897567749, -1.859290599822998, 0.42024119551909, -0.5193942967107517, -1.513404460276592),
    @(7, -6.19073486328125, -5.0089111328125, -5.00853443145752, 0.1340822670293074, -0.154626976635499, -1.524474680423737),
    @(8, -16.53120613098145, -4.607099056243896, 0.276792049407959, -0.4547065528015019, -0.6670010999097716, 0.1219344166881118),
    @(9, -4.55709171295166, -19.43513488769531, -7.060164451599121, -0.1632935243769215, -1.267501022350991, -1.733805928316158),
    @(10, -1.290600776672363, -3.404523849487305, 9.77668571472168, 1.678417857420652, 1.337361679238797, -3.577607195256125),
    @(11, -2.725464344024658, -15.71722602844238, -1.773126602172852, -1.108878275600511, -0.9978605397676845, 0.5690611009880358),
    @(12, -0.2208814620971679, -29.22537994384766, -4.326999664306641, -1.319620092036361, 2.702669687190279, 3.955988294997444),
    @(13, -30.18045997619629, -3.106259346008301, -5.336086273193359, -4.884191094818762, 2.429627870099019, 5.183373525991279),
    @(14, -28.01585388183594, -10.60748767852783, -5.628327369689941, 0.9754154722569472, -2.997867519572641, -0.7411238864316898),
    @(15, 3.033831596374512, -20.5565071105957, -19.63080787658692, 2.564180643376662, -3.320737330590262, -3.21195284693929),
    @(16, 46.15228271484375, 6.39080286026001, 16.86857604980469, 4.649666024466685, -1.874780806444464, -0.4834878709704054),
    @(17, -12.0808572769165, 1.62255322933197, -1.056400775909424, 1.0025165415669, 0.4837598439495385, -3.280407503113078),
    @(18, -11.0537633895874, -42.16847610473633, 6.546759128570557, -4.006519524727852, 1.548193707304488, 2.943372968900025),
    @(19, -1.753406047821045, 1.232700347900391, -26.96573448181152, -9.41962250814602, 4.90418978969931, 3.12689200195215),
    @(20, -24.26501083374023, -16.1025276184082, 0.0206184387207031, -1.074812318308926, -4.687652131258382, 0.3980861830256872),
    @(21, 6.58524227142334, -32.96889877319336, 19.44017219543457, 2.882893306218964, -2.573564957764209, -5.314060277289736),
    @(22, 8.066184997558594, 4.324060440063477, 3.93980073928833, 6.499037447622277, 4.13717439619163, -5.053133794816869),
    @(23, 6.809474468231201, -5.948239803314209, 5.694517135620117, 0.6406787679356816, 1.996182547282352, -0.02198917330320072),
    @(24, -36.60982131958008, -23.81085586547852, -32.70140838623047, -5.474647755340023, -0.02923721226595877, 3.930300022585936),
    @(25, 23.85481452941895, 0.3514032363891601, -14.82476615905762, -3.350842246564774, 7.508822642393016, 1.570750523421729),
    @(26, -12.4519100189209, -15.55856132507324, -9.169242858886721, -1.041599551499903, 6.251372620210817, -1.606606071278196),
    @(27, 23.85086059570312, -22.51275253295898, 13.36398124694824, 4.526515125470683, -2.61241545313501, -4.543736482070679),
    @(28, -5.659902095794678, 5.291859149932861, -6.011741161346436, 5.028566198833909, -0.9918887261616698, -6.872307486453308),
    @(29, 32.19543075561523, -80.94894409179688, 22.37196731567383, 0.376901606381959, 1.241255387916403, -0.1712141006680019),
    @(30, -18.56385040283203, 7.148910522460938, -19.82358169555664, -8.190219378067249, -11.05336972818532, 0.1851687067646206),
    @(31, -14.20774078369141, -7.775349617004394, -14.13762283325195, 0.4136780055901301, -15.04960784669652, 3.795090385412704)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
    $ws.Cells.Item($r, 2).Value = "walkingToRunning"
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
}
